# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# In the "Hoja1" sheet, the table starting at row 15 lists, for each
# "Periodo Mora" (year), its corresponding "Valor Mora" in column F:
#   Row 16 -> 2007
#   Row 17 -> 2006
#   Row 18 -> 2005
#   Row 19 -> 2004
#   Row 20 -> 2003
#   Row 21 -> 2002
#
# The database refresh swapped the "Valor Mora" amounts recorded for the
# 2007 period (row 16) and the 2002 period (row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F16").Value = 33125
$ws.Range("F21").Value = 24292
